# Vikram_biodata.docx edit:
#   The father's surname "Panday" (in "Father ... Devanand <tabs> Panday")
#   is corrected to "Pandey" by changing the 5th letter from "a" to "e".
#   In the canonical OOXML this shows up as the single run <w:t>Panday</w:t>
#   being split into three runs - "Pand" (the original run, now holding
#   only the unchanged prefix), "e" (brand new run) and "y" (brand new
#   run) - all three sharing the same run formatting.
#
# There are three occurrences of "Panday" in this document (Vikash Panday,
# Suraj Prakas Panday, and the father Devanand Panday); we disambiguate by
# anchoring on "Devanand", which immediately precedes the target text.

$d = $word.ActiveDocument
$full = $d.Content.Text

$anchor = "Devanand"
$anchorIdx = $full.IndexOf($anchor)
if ($anchorIdx -lt 0) { throw "Could not locate anchor text 'Devanand'" }

$searchFrom = $anchorIdx + $anchor.Length
$targetIdx = $full.IndexOf("Panday", $searchFrom)
if ($targetIdx -lt 0) { throw "Could not locate 'Panday' after 'Devanand'" }

# "Panday" character offsets relative to $targetIdx: P=0 a=1 n=2 d=3 a=4 y=5
$aPos = $targetIdx + 4   # the "a" that must become "e"
$yPos = $targetIdx + 5   # the trailing "y"

$origColor = $d.Range($aPos, $aPos + 1).Font.Color

# Step 1: temporarily give the "a" and the "y" two different (and different
# from the surrounding text) font colors. This keeps the text engine from
# silently re-merging them into their identically-formatted neighbouring
# run once we touch them, so each ends up materialised as its own <w:r>.
$rA = $d.Range($aPos, $aPos + 1)
$rA.Font.Color = 255
$rY = $d.Range($yPos, $yPos + 1)
$rY.Font.Color = 65280

# Step 2: rewrite the "a" run's text to "e" (Panday -> Pandey) and force a
# content churn through the "y" run as well (via a throwaway placeholder
# character, since setting it to its own existing text is a no-op that the
# engine would not treat as a fresh edit), then restore both runs' font
# color back to the original value. Each of these two runs now persists as
# an independent, freshly-minted run alongside the untouched "Pand" run.
$rA2 = $d.Range($aPos, $aPos + 1)
$rA2.Text = "e"
$rA3 = $d.Range($aPos, $aPos + 1)
$rA3.Font.Color = $origColor

$rY2 = $d.Range($yPos, $yPos + 1)
$rY2.Text = "Q"
$rY2b = $d.Range($yPos, $yPos + 1)
$rY2b.Text = "y"
$rY3 = $d.Range($yPos, $yPos + 1)
$rY3.Font.Color = $origColor

$result = $d.Range($targetIdx, $targetIdx + 6)
Write-Output "Panday -> $($result.Text)"
